# Synchronize the results table: rotate the per-matrix benchmark blocks
# (rows 2-21, grouped in pairs of 2 rows per matrix) and refresh the
# "Tempo" (time) column with newly measured values, starting at the
# fifth row-group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for columns B (Matriz), D (Autovalor), E (Iteracoes),
# F (Tempo), G (Ordem), H (Campo), I (Simetria) for rows 2..21.
$data = @(
    @{ Row=2;  B="c130.mtx";     D=2.366980478285793;  E=113; F=0.004698038101196289; G=130;  H="real";    I="general"   },
    @{ Row=3;  B="c130.mtx";     D=2.370733342517338;  E=39;  F=0.001653671264648438; G=130;  H="real";    I="general"   },
    @{ Row=4;  B="bcspwr04.mtx"; D=10.11425564616074;  E=30;  F=0.001617193222045898; G=274;  H="pattern"; I="symmetric" },
    @{ Row=5;  B="bcspwr04.mtx"; D=10.11401361716994;  E=15;  F=0.0008158683776855469; G=274; H="pattern"; I="symmetric" },
    @{ Row=6;  B="bcspwr02.mtx"; D=4.308258605005528;  E=40;  F=0.001262664794921875; G=49;   H="pattern"; I="symmetric" },
    @{ Row=7;  B="bcspwr02.mtx"; D=4.308461254391954;  E=24;  F=0.0008020401000976562; G=49;  H="pattern"; I="symmetric" },
    @{ Row=8;  B="bcspwr03.mtx"; D=5.104674213788415;  E=35;  F=0.001567840576171875; G=118;  H="pattern"; I="symmetric" },
    @{ Row=9;  B="bcspwr03.mtx"; D=5.104142633099601;  E=11;  F=0.000446319580078125; G=118;  H="pattern"; I="symmetric" },
    @{ Row=10; B="38_bus.mtx";   D=30115.82629984464;  E=22;  F=0.01273417472839355;  G=1138; H="real";    I="symmetric" },
    @{ Row=11; B="38_bus.mtx";   D=30114.18326043027;  E=14;  F=0.009031057357788086; G=1138; H="real";    I="symmetric" },
    @{ Row=12; B="685_bus.mtx";  D=26186.48629082991;  E=4;   F=0.001051425933837891; G=685;  H="real";    I="symmetric" },
    @{ Row=13; B="685_bus.mtx";  D=26186.48629088882;  E=4;   F=0.0008580684661865234; G=685; H="real";    I="symmetric" },
    @{ Row=14; B="h292.mtx";     D=9.151774331566605;  E=45;  F=0.002381324768066406; G=292;  H="pattern"; I="symmetric" },
    @{ Row=15; B="h292.mtx";     D=9.152249584058874;  E=23;  F=0.001132488250732422; G=292;  H="pattern"; I="symmetric" },
    @{ Row=16; B="494_bus.mtx";  D=2220.967903929713;  E=4;   F=0.0006465911865234375; G=494; H="real";    I="symmetric" },
    @{ Row=17; B="494_bus.mtx";  D=2220.970984428201;  E=4;   F=0.0004987716674804688; G=494; H="real";    I="symmetric" },
    @{ Row=18; B="h85.mtx";      D=6.721075798530548;  E=27;  F=0.001138925552368164; G=85;   H="pattern"; I="symmetric" },
    @{ Row=19; B="h85.mtx";      D=6.721128139935032;  E=18;  F=0.0008838176727294922; G=85;  H="pattern"; I="symmetric" },
    @{ Row=20; B="662_bus.mtx";  D=1422.994613068716;  E=4;   F=0.0009610652923583984; G=662;  H="real";    I="symmetric" },
    @{ Row=21; B="662_bus.mtx";  D=1422.994613069662;  E=4;   F=0.0008969306945800781; G=662;  H="real";    I="symmetric" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
    $ws.Range("H$r").Value = $item.H
    $ws.Range("I$r").Value = $item.I
}
